# Applies the StructureDefinition-claim-category.xlsx update:
#  - Metadata sheet: bump Version, Date, Publisher, replace duplicate
#    "Contact" rows with "Jurisdiction"/"United States of America", and
#    remove the now-redundant extra row so the table is 20 rows again.
#  - Elements sheet: update the root Extension row's Short/Definition
#    text to describe "Claim Category" instead of the generic
#    "Extension"/"An Extension" placeholders.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value  = "6.0.0"
$meta.Range("B8").Value  = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value  = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" row; delete it so everything below
# shifts up by one (table goes from 21 rows to 20 rows).
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Claim Category"
$elements.Range("L2").Value = "Standard HIPAA code for the category of the claim status"
